# Automatic update of files.
# The data rows (2-5) are cyclically rotated up by one: the record that was
# in row 2 moves down to row 5, and the records in rows 3, 4, 5 each move up
# by one row (to 2, 3, 4 respectively). Apply the change by writing the
# already-rotated values directly into each cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the original values of each affected column for rows 2-5 before
# overwriting anything.
$cols = @("A","B","D","E","F","G","H","Q","R","S","Z","AB","AW","AX")

$orig = @{}
foreach ($col in $cols) {
    $orig[$col] = @{}
    for ($r = 2; $r -le 5; $r++) {
        $orig[$col][$r] = $ws.Range("$col$r").Value2
    }
}

# new row 2 <- old row 3
# new row 3 <- old row 4
# new row 4 <- old row 5
# new row 5 <- old row 2
$srcRow = @{ 2 = 3; 3 = 4; 4 = 5; 5 = 2 }

foreach ($col in $cols) {
    foreach ($destRow in 2..5) {
        $source = $srcRow[$destRow]
        $ws.Range("$col$destRow").Value2 = $orig[$col][$source]
    }
}
